$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2262.077
$ws.Range("J18").Value = 2401
$ws.Range("L18").Value = 2401
$ws.Range("N18").Value = -2969

$ws.Range("H51").Value = 5838.4443
$ws.Range("J51").Value = 5945.375
$ws.Range("L51").Value = 5945.375
$ws.Range("N51").Value = -6913.375

$ws.Range("H80").Value = 144.11111
$ws.Range("I80").Value = 46
$ws.Range("J80").Value = 222.6
$ws.Range("K80").Value = 138
$ws.Range("L80").Value = 667.8
$ws.Range("M80").Value = 860
$ws.Range("N80").Value = -2663.8

$ws.Range("H83").Value = 144.11111
$ws.Range("I83").Value = 46
$ws.Range("J83").Value = 222.6
$ws.Range("K83").Value = 414
$ws.Range("L83").Value = 2003.4
$ws.Range("M83").Value = 4578
$ws.Range("N83").Value = -11987.4

$ws.Range("H88").Value = 336080.84
$ws.Range("I88").Value = 402397
$ws.Range("K88").Value = 402397
$ws.Range("M88").Value = -401991

$ws.Range("H91").Value = 336080.84
$ws.Range("I91").Value = 402397
$ws.Range("K91").Value = 402397
$ws.Range("M91").Value = -400993

$ws.Range("H98").Value = 3259.2856
$ws.Range("I98").Value = 2952.5
$ws.Range("K98").Value = 2952.5
$ws.Range("M98").Value = -1454.5

$ws.Range("H100").Value = 996.6667
$ws.Range("I100").Value = 990
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 990
$ws.Range("L100").Value = 1000
$ws.Range("M100").Value = -449
$ws.Range("N100").Value = -2082

$ws.Range("H122").Value = 3259.2856
$ws.Range("I122").Value = 2952.5
$ws.Range("K122").Value = 8857.5
$ws.Range("M122").Value = -6407.5

$ws.Range("H129").Value = 4816.3335
$ws.Range("I129").Value = 5176.8
$ws.Range("J129").Value = 3014
$ws.Range("K129").Value = 15530.4
$ws.Range("L129").Value = 9042
$ws.Range("M129").Value = -10530.4
$ws.Range("N129").Value = -19042

$ws.Range("H132").Value = 3034.1904
$ws.Range("I132").Value = 2761.0625
$ws.Range("J132").Value = 3908.2
$ws.Range("K132").Value = 8283.1875
$ws.Range("L132").Value = 11724.6
$ws.Range("M132").Value = -5753.1875
$ws.Range("N132").Value = -16784.6

$ws.Range("H137").Value = 4182.033
$ws.Range("J137").Value = 5755.091
$ws.Range("L137").Value = 17265.273
$ws.Range("N137").Value = -22365.273

$ws.Range("H141").Value = 5769.143
$ws.Range("I141").Value = 4770.5
$ws.Range("J141").Value = 7100.6665
$ws.Range("K141").Value = 14311.5
$ws.Range("L141").Value = 21301.9995
$ws.Range("M141").Value = -9131.5
$ws.Range("N141").Value = -31661.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 5250
$ws.Range("I36").Value = 5250
$ws.Range("K36").Value = 5250
$ws.Range("M36").Value = -4904

$ws.Range("H61").Value = 4988.231
$ws.Range("I61").Value = 4443.5
$ws.Range("J61").Value = 5859.8
$ws.Range("K61").Value = 4443.5
$ws.Range("L61").Value = 5859.8
$ws.Range("M61").Value = -4231.5
$ws.Range("N61").Value = -6283.8

$ws.Range("H132").Value = 2588.6428
$ws.Range("I132").Value = 2061.4546
$ws.Range("J132").Value = 4521.6665
$ws.Range("K132").Value = 6184.3638
$ws.Range("L132").Value = 13564.9995
$ws.Range("M132").Value = -3654.3638
$ws.Range("N132").Value = -18624.9995

$ws.Range("H136").Value = 4988.231
$ws.Range("I136").Value = 4443.5
$ws.Range("J136").Value = 5859.8
$ws.Range("K136").Value = 13330.5
$ws.Range("L136").Value = 17579.4
$ws.Range("M136").Value = -10780.5
$ws.Range("N136").Value = -22679.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1230.5
$ws.Range("I20").Value = 1256.2222
$ws.Range("K20").Value = 1256.2222
$ws.Range("M20").Value = -1009.2222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3124.5
$ws.Range("I31").Value = 2094.375
$ws.Range("J31").Value = 5184.75
$ws.Range("K31").Value = 2094.375
$ws.Range("L31").Value = 5184.75
$ws.Range("M31").Value = -1799.375
$ws.Range("N31").Value = -5774.75

$ws.Range("H32").Value = 832.3333
$ws.Range("I32").Value = 832.3333
$ws.Range("K32").Value = 832.3333
$ws.Range("M32").Value = -516.3333

$ws.Range("H34").Value = 3124.5
$ws.Range("I34").Value = 2094.375
$ws.Range("J34").Value = 5184.75
$ws.Range("K34").Value = 2094.375
$ws.Range("L34").Value = 5184.75
$ws.Range("M34").Value = -1892.375
$ws.Range("N34").Value = -5588.75

$ws.Range("H58").Value = 3167.5789
$ws.Range("I58").Value = 2967.4666
$ws.Range("J58").Value = 3918
$ws.Range("K58").Value = 2967.4666
$ws.Range("L58").Value = 3918
$ws.Range("M58").Value = -2764.4666
$ws.Range("N58").Value = -4324

$ws.Range("H99").Value = 138199.2
$ws.Range("I99").Value = 90499.5
$ws.Range("K99").Value = 90499.5
$ws.Range("M99").Value = -89001.5

$ws.Range("H107").Value = 1305.75
$ws.Range("J107").Value = 1610.625
$ws.Range("L107").Value = 1610.625
$ws.Range("N107").Value = -5450.625

$ws.Range("H122").Value = 2003.6
$ws.Range("J122").Value = 1717.5
$ws.Range("L122").Value = 5152.5
$ws.Range("N122").Value = -10052.5

$ws.Range("H126").Value = 138199.2
$ws.Range("I126").Value = 90499.5
$ws.Range("K126").Value = 271498.5
$ws.Range("M126").Value = -269028.5

$ws.Range("H134").Value = 3502.5715
$ws.Range("I134").Value = 3502.5715
$ws.Range("K134").Value = 10507.7145
$ws.Range("M134").Value = -7972.7145

$ws.Range("H136").Value = 3167.5789
$ws.Range("I136").Value = 2967.4666
$ws.Range("J136").Value = 3918
$ws.Range("K136").Value = 8902.399800000001
$ws.Range("L136").Value = 11754
$ws.Range("M136").Value = -6352.399800000001
$ws.Range("N136").Value = -16854

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 1388.1818
$ws.Range("J9").Value = 1388.1818
$ws.Range("L9").Value = 4164.5454
$ws.Range("N9").Value = -4612.5454

$ws.Range("H68").Value = 2495
$ws.Range("I68").Value = 2966.875
$ws.Range("J68").Value = 1865.8334
$ws.Range("K68").Value = 8900.625
$ws.Range("L68").Value = 5597.5002
$ws.Range("M68").Value = -8089.625
$ws.Range("N68").Value = -7219.5002

$ws.Range("H71").Value = 2495
$ws.Range("I71").Value = 2966.875
$ws.Range("J71").Value = 1865.8334
$ws.Range("K71").Value = 26701.875
$ws.Range("L71").Value = 16792.5006
$ws.Range("M71").Value = -22645.875
$ws.Range("N71").Value = -24904.5006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 11000
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

$ws.Range("H70").Value = 6637.4614
$ws.Range("J70").Value = 7158.9
$ws.Range("L70").Value = 7158.9
$ws.Range("N70").Value = -7698.9

$ws.Range("H73").Value = 6637.4614
$ws.Range("J73").Value = 7158.9
$ws.Range("L73").Value = 7158.9
$ws.Range("N73").Value = -9030.9

$ws.Range("H97").Value = 1127
$ws.Range("I97").Value = 1127
$ws.Range("K97").Value = 1127
$ws.Range("M97").Value = -631

$ws.Range("H122").Value = 2167.6667
$ws.Range("I122").Value = 2167.6667
$ws.Range("K122").Value = 6503.000100000001
$ws.Range("M122").Value = -4053.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3832.4546
$ws.Range("I22").Value = 2787.7144
$ws.Range("J22").Value = 4113.731
$ws.Range("K22").Value = 2787.7144
$ws.Range("L22").Value = 4113.731
$ws.Range("M22").Value = -2492.7144
$ws.Range("N22").Value = -4703.731

$ws.Range("H27").Value = 3832.4546
$ws.Range("I27").Value = 2787.7144
$ws.Range("J27").Value = 4113.731
$ws.Range("K27").Value = 2787.7144
$ws.Range("L27").Value = 4113.731
$ws.Range("M27").Value = -2680.7144
$ws.Range("N27").Value = -4327.731

$ws.Range("H46").Value = 2069.818
$ws.Range("I46").Value = 969.5
$ws.Range("J46").Value = 2698.5715
$ws.Range("K46").Value = 969.5
$ws.Range("L46").Value = 2698.5715
$ws.Range("M46").Value = -781.5
$ws.Range("N46").Value = -3074.5715

$ws.Range("H55").Value = 351.4737
$ws.Range("I55").Value = 382.30768
$ws.Range("J55").Value = 284.66666
$ws.Range("K55").Value = 382.30768
$ws.Range("L55").Value = 284.66666
$ws.Range("M55").Value = -209.30768
$ws.Range("N55").Value = -630.66666

$ws.Range("H82").Value = 2099.4736
$ws.Range("I82").Value = 2165.9285
$ws.Range("J82").Value = 1913.4
$ws.Range("K82").Value = 2165.9285
$ws.Range("L82").Value = 1913.4
$ws.Range("M82").Value = -1804.9285
$ws.Range("N82").Value = -2635.4

$ws.Range("H85").Value = 2099.4736
$ws.Range("I85").Value = 2165.9285
$ws.Range("J85").Value = 1913.4
$ws.Range("K85").Value = 2165.9285
$ws.Range("L85").Value = 1913.4
$ws.Range("M85").Value = -917.9285
$ws.Range("N85").Value = -4409.4

$ws.Range("H122").Value = 4852.6924
$ws.Range("I122").Value = 4208.8
$ws.Range("K122").Value = 12626.4
$ws.Range("M122").Value = -10176.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 9000
$ws.Range("J47").Value = 9000
$ws.Range("L47").Value = 9000
$ws.Range("N47").Value = -10144

$ws.Range("H125").Value = 67500
$ws.Range("J125").Value = 67500
$ws.Range("L125").Value = 67500
$ws.Range("N125").Value = -77340

$ws.Range("H126").Value = 2667.2354
$ws.Range("I126").Value = 2929.6086
$ws.Range("J126").Value = 2118.6365
$ws.Range("K126").Value = 8788.825800000001
$ws.Range("L126").Value = 6355.9095
$ws.Range("M126").Value = -6318.825800000001
$ws.Range("N126").Value = -11295.9095

$ws.Range("H132").Value = 3063.8103
$ws.Range("I132").Value = 2965.5273
$ws.Range("K132").Value = 8896.581900000001
$ws.Range("M132").Value = -6366.581900000001
